$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new "2023" column (K) to the table, mirroring the formatting of the
# preceding "2022" column (J): header year, and the three data rows
# (Gel / Women / Men).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 613.5
$ws.Range("K5").Value = 395.8
$ws.Range("K6").Value = 737.5
